$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Exp date"
$ws.Range("B1").Value = "Data to be processed (path)"
$ws.Range("C1").Value = "Path of the pillar array info"
$ws.Range("D1").Value = "Saving path"
$ws.Range("E1").Value = "Channel height (um)"
$ws.Range("F1").Value = "Channel width (um)"
$ws.Range("G1").Value = "Flow rate (nL/s)"
$ws.Range("H1").Value = "Initial velocity (m/s)"
$ws.Range("I1").Value = "Calibration (um/pixel)"

# ---- Row 2 (20210413-Actin) ----
$ws.Range("B2").Value = "F:\PhD, PMMH, ESPCI\Processing\20210413-Actin\results\Shapes_1"
$ws.Range("C2").Value = "F:\PhD, PMMH, ESPCI\Processing\20210413-Actin\results\circlesforPAs1_S10.mat"
$ws.Range("D2").Value = "F:\PhD, PMMH, ESPCI\Processing\20210413-Actin\results\Figures"
$ws.Range("E2").Value = 48
$ws.Range("F2").Value = 400
$ws.Range("G2").Value = 0.5
$ws.Range("H2").Formula = "=G2/F2/E2"
$ws.Range("I2").Value = 0.1

# ---- Row 3 (20210430-Actin, Shapes_1) ----
$ws.Range("B3").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Shapes_1"
$ws.Range("C3").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\circlesforPAs2_S10.mat"
$ws.Range("D3").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Figures"
$ws.Range("E3").Value = 55
$ws.Range("F3").Value = 400
$ws.Range("G3").Value = 0.5

# ---- Row 4 (20210430-Actin, Shapes_2) ----
$ws.Range("B4").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Shapes_2"
$ws.Range("C4").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\circlesforPAs2_S10.mat"
$ws.Range("D4").Value = "F:\PhD, PMMH, ESPCI\Processing\20210430-Actin\results\Figures"
$ws.Range("E4").Value = 55
$ws.Range("F4").Value = 400
$ws.Range("G4").Value = 1

# Shared formula for H3:H4
$ws.Range("H3:H4").Formula = "=G3/F3/E3"

$ws.Range("I3").Value = 0.1
$ws.Range("I4").Value = 0.1

# ---- Exp date column with date number format (re-use one style via copy/paste of formats) ----
$ws.Range("A2").Value = 44299
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$ws.Range("A3").Value = 44316
$ws.Range("A4").Value = 44316
$excel.CutCopyMode = $false

# ---- Column widths ----
# NOTE: this host's ColumnWidth setter quantizes to the nearest 1/6 of a
# character *after* adding a fixed 5/6-character padding (i.e. stored
# width = round(6*input)/6 + 5/6). We pre-compensate the input by that
# same 5/6 offset so the value actually written to the sheet lands on the
# closest achievable 1/6-quantum to the real target width.
$ws.Columns.Item(1).ColumnWidth = 8.833333333333334   # target 9.7109375
$ws.Columns.Item(2).ColumnWidth = 63.666666666666664  # target 64.5703125
$ws.Columns.Item(3).ColumnWidth = 66.66666666666667   # target 67.42578125
$ws.Columns.Item(4).ColumnWidth = 59.166666666666664  # target 60
$ws.Columns.Item(5).ColumnWidth = 18.666666666666668  # target 19.42578125
$ws.Columns.Item(6).ColumnWidth = 18.0                # target 18.85546875
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666  # target 15
$ws.Columns.Item(8).ColumnWidth = 18.333333333333332  # target 19.140625
$ws.Columns.Item(9).ColumnWidth = 20.166666666666668  # target 21

# ---- View / selection state ----
[void]$ws.Range("E18").Select()
